# Ajeitando o formato dos graficos de comparação
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "C2" "25°"
Set-TextValue "D3" "58%"
Set-TextValue "E3" "78%"
Set-TextValue "B4" "35°"
Set-TextValue "C4" "25°"
Set-TextValue "B8" "34°"
Set-TextValue "D9" "72%"
Set-TextValue "E9" "90%"
